# Applies the "Wind turbine onshore" process update to VT_Model_ALL_V01.xlsx
$wb = $excel.ActiveWorkbook

$wsProcesses = $wb.Worksheets.Item("SEC_Processes")
$wsPP        = $wb.Worksheets.Item("PP")

# --- SEC_Processes: new technology row (row 10) ---------------------------
$wsProcesses.Activate()
$wsProcesses.Range("B10").Value = "ELE"
$wsProcesses.Range("E10").Value = "Wind turbine onshore"
$wsProcesses.Range("D10").Value = "ELE_EX_WIND_TURBINE"
$wsProcesses.Range("F10").Value = "PJ"
$wsProcesses.Range("G10").Value = "GWe"
$wsProcesses.Range("H10").Value = "DAYNITE"

$wsProcesses.Range("D11").Select() | Out-Null

# --- PP sheet: populate row 9 with the new technology's parameters --------
# H9 carries the same format as H8 (no shading/border) rather than the
# thick-bottom row format it currently has, so copy that formatting over
# before writing the value.
$wsPP.Range("H8").Copy()
$wsPP.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPP.Range("B9").Formula = "=SEC_Processes!D10"
$wsPP.Range("C9").Formula = "=SEC_Processes!E10"
$wsPP.Range("D9").Formula = "=SEC_Comm!C9"
$wsPP.Range("E9").Formula = "=SEC_Comm!C8"
$wsPP.Range("F9").Value = 1.345
$wsPP.Range("G9").Value = 1
$wsPP.Range("H9").Value = 31.536
$wsPP.Range("I9").Value = 0.33
$wsPP.Range("J9").Value = 1
$wsPP.Range("K9").Value = 0

# --- PP sheet: derived calculation rows ------------------------------------
$wsPP.Range("C12").Value = "Maximum output"
$wsPP.Range("D12").Formula = "=F9*H9"
$wsPP.Range("D12").ClearFormats()
$wsPP.Range("E12").Value = "PJ"

$wsPP.Range("C13").Value = "Limitet output"
$wsPP.Range("D13").Formula = "=D12*I9"
$wsPP.Range("D13").ClearFormats()
$wsPP.Range("E13").Value = "PJ"

# --- Activate PP sheet and set the new selection / zoom --------------------
$wsPP.Activate()
$wsPP.Range("D17").Select() | Out-Null
$excel.ActiveWindow.Zoom = 170
